# Insert a new worksheet "03_Factorization" between "02_Law of Indices" and
# "20_Properties of Circles", populate it with its question-code data, and
# tidy up the selections on the sheets that shift as a result.

$wb = $excel.ActiveWorkbook

# The sheet that currently sits right after where the new sheet must go.
$ws2 = $wb.Worksheets.Item("02_Law of Indices")

# Reset sheet 2's selection to its full data range (loses tabSelected once
# the new sheet becomes active below).
[void]$ws2.Range("A1:B15").Select()

# Insert the new worksheet immediately after "02_Law of Indices".
$newSheet = $wb.Worksheets.Add($null, $ws2)
$newSheet.Name = "03_Factorization"

# Populate the new sheet with its question-code table.
$data = @(
    @(3, "Question Code"),
    @(1, "DSE17PII_Q01"),
    @(2, "DSE23PII_Q04"),
    @(3, "DSE24PII_Q01"),
    @(4, "DSE13PII_Q03"),
    @(5, "DSE20PII_Q04"),
    @(6, "DSE18PII_Q03"),
    @(7, "DSE14PII_Q02"),
    @(8, "DSE22PII_Q01"),
    @(9, "DSE12PII_Q02"),
    @(10, "DSE16PII_Q03"),
    @(11, "DSEPPPII_Q03"),
    @(12, "DSESPPII_Q03")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
}

# Match the sheet's final selection/active-cell state.
[void]$newSheet.Range("A14:B15").Select()
